# This script re-arranges (permutes) the data held in rows 17-30 of the
# "Artfynd" worksheet. Every cell in a given source row moves, as a whole
# row, to a (generally different) target row. One row (21) stays in place.
#
# Mapping below reads as: target row = source row (i.e. the content that
# used to live in the source row ends up, after the edit, in the target
# row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    17 = 20
    18 = 23
    19 = 26
    20 = 25
    21 = 21
    22 = 30
    23 = 19
    24 = 29
    25 = 24
    26 = 18
    27 = 17
    28 = 22
    29 = 28
    30 = 27
}

$firstCol = 1   # A
$lastCol = 51   # AY

# Helper: assign a value to a cell without letting Excel's "looks like a
# date" auto-detection silently convert a literal text value (e.g.
# "2023-09-17") into a date serial number. We briefly force a text number
# format for the assignment, then restore the cell to the default "Normal"
# style so no stray style index is left behind.
function Set-CellValueSafe {
    param($Cell, $Val)

    if ($Val -eq $null) {
        $Cell.Value = ""
        return
    }
    if ($Val -is [string] -and $Val -match '^\d{4}-\d{1,2}-\d{1,2}$') {
        $Cell.NumberFormat = "@"
        $Cell.Value = $Val
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Val
    }
}

# Step 1: snapshot every source row's values (all columns) before writing
# anything back, because several rows are both a source for one target and
# a target that receives data from another row.
$snapshot = @{}
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $rowValues = New-Object 'object[]' ($lastCol - $firstCol + 1)
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowValues[$col - $firstCol] = $ws.Cells.Item($sourceRow, $col).Value()
    }
    $snapshot[$targetRow] = $rowValues
}

# Step 2: write the snapshotted values into their target rows.
foreach ($targetRow in $rowMap.Keys) {
    $rowValues = $snapshot[$targetRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $v = $rowValues[$col - $firstCol]
        Set-CellValueSafe $ws.Cells.Item($targetRow, $col) $v
    }
}
